# Update Test Data used in data driven tests (SeleniumWithJava TestData1.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test data values on row 8
$ws.Range("A8").Value = "TestAutomation_POC5"
$ws.Range("B8").Value = "TestAutomation_POC5"
$ws.Range("C8").Value = "Facility_POC5"
$ws.Range("D8").Value = "Facility_POC5"
$ws.Range("E8").Value = "Pharmacy_POC5"
$ws.Range("F8").Value = "Pharmacy_POC5"
$ws.Range("H8").Value = "Alignment Project _POC5"

# Update the active cell selection on the sheet
$ws.Activate()
$ws.Range("H19").Select()
